# Correct petroleum's BAU dispatch priority: it should guarantee dispatch
# (~10%) like the peakers do, so flip its priority order from 2 (merit-order)
# to 1 (must-run / guaranteed-dispatch), matching the same treatment as
# natural gas peakers - without actually flagging it as a peaker.
#
# Column B holds the literal priority value for 2015; columns C:AK are
# shared formulas "=$B<row>" that copy it across every subsequent year, so
# updating B11 alone ripples the value through the rest of the row.

$wb = $excel.ActiveWorkbook

$wsBDPbES = $wb.Worksheets.Item("BDPbES")

# petroleum is row 11 on the BDPbES sheet
$wsBDPbES.Range("B11").Value = 1

# The workbook was last saved with the BDPbES sheet active/selected (and
# About no longer the active tab), with the selection left on H17.
$wsBDPbES.Activate()
$wsBDPbES.Range("H17").Select()
